$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 75

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-19"
$ws.Cells.Item($row, 2).Value = "15:17:25"
$ws.Cells.Item($row, 3).Value = "Friday"
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "02"
$ws.Cells.Item($row, 5).Value = 138107
$ws.Cells.Item($row, 6).Value = 140427
$ws.Cells.Item($row, 7).Value = 171324
$ws.Cells.Item($row, 8).Value = 148859
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 121873
$ws.Cells.Item($row, 11).Value = 223503
$ws.Cells.Item($row, 12).Value = 254914
$ws.Cells.Item($row, 13).Value = 185258
$ws.Cells.Item($row, 14).Value = 110357
$ws.Cells.Item($row, 15).Value = 41346
$ws.Cells.Item($row, 16).Value = 30906
$ws.Cells.Item($row, 17).Value = 73555
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42556
$ws.Cells.Item($row, 20).Value = -1
